# Apply updated odds/stats values to row 5 of the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 3.1
$ws.Range("I5").Value = 2.2
$ws.Range("J5").Value = 3.6
$ws.Range("U5").Value = 1.62
$ws.Range("V5").Value = 2.2
$ws.Range("Y5").Value = 12
$ws.Range("AC5").Value = 12
$ws.Range("AG5").Value = 151
$ws.Range("AH5").Value = 9.5
$ws.Range("AO5").Value = 17
$ws.Range("AP5").Value = 23
$ws.Range("AS5").Value = 151
$ws.Range("AV5").Value = 51
$ws.Range("AW5").Value = 4.33
